{"js": "// Update the worksheet date header and every two-digit-by-two-digit\n// multiplication problem/answer cell to the new day's values.\nconst replacements = [\n  [\"2025-01-22 Wednesday\", \"2025-01-23 Thursday\"],\n  [\"81\u00d790=7290\", \"42\u00d792=3864\"],\n  [\"90\u00d712=1080\", \"23\u00d712=276\"],\n  [\"87\u00d793=8091\", \"28\u00d768=1904\"],\n  [\"98\u00d783=8134\", \"20\u00d733=660\"],\n  [\"21\u00d792=1932\", \"99\u00d731=3069\"],\n  [\"15\u00d746=690\", \"18\u00d775=1350\"],\n  [\"83\u00d763=5229\", \"61\u00d757=3477\"],\n  [\"50\u00d732=1600\", \"65\u00d785=5525\"],\n  [\"65\u00d761=3965\", \"48\u00d783=3984\"],\n  [\"95\u00d743=4085\", \"78\u00d797=7566\"],\n  [\"77\u00d735=2695\", \"77\u00d749=3773\"],\n  [\"46\u00d711=506\", \"74\u00d765=4810\"],\n  [\"36\u00d751=1836\", \"64\u00d757=3648\"],\n  [\"37\u00d766=2442\", \"77\u00d776=5852\"],\n  [\"94\u00d729=2726\", \"67\u00d761=4087\"],\n  [\"67\u00d715=1005\", \"17\u00d719=323\"],\n  [\"89\u00d762=5518\", \"92\u00d720=1840\"],\n  [\"63\u00d754=3402\", \"59\u00d752=3068\"],\n  [\"59\u00d720=1180\", \"25\u00d713=325\"],\n  [\"87\u00d740=3480\", \"22\u00d724=528\"],\n  [\"32\u00d787=2784\", \"92\u00d725=2300\"],\n  [\"60\u00d734=2040\", \"54\u00d772=3888\"],\n  [\"69\u00d794=6486\", \"16\u00d780=1280\"],\n  [\"78\u00d719=1482\", \"20\u00d782=1640\"],\n  [\"84\u00d783=6972\", \"65\u00d734=2210\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and every two-digit-by-two-digit\n# multiplication problem/answer cell to the new day's values.\n$pairs = @(\n    @(\"2025-01-22 Wednesday\", \"2025-01-23 Thursday\"),\n    @(\"81\u00d790=7290\", \"42\u00d792=3864\"),\n    @(\"90\u00d712=1080\", \"23\u00d712=276\"),\n    @(\"87\u00d793=8091\", \"28\u00d768=1904\"),\n    @(\"98\u00d783=8134\", \"20\u00d733=660\"),\n    @(\"21\u00d792=1932\", \"99\u00d731=3069\"),\n    @(\"15\u00d746=690\", \"18\u00d775=1350\"),\n    @(\"83\u00d763=5229\", \"61\u00d757=3477\"),\n    @(\"50\u00d732=1600\", \"65\u00d785=5525\"),\n    @(\"65\u00d761=3965\", \"48\u00d783=3984\"),\n    @(\"95\u00d743=4085\", \"78\u00d797=7566\"),\n    @(\"77\u00d735=2695\", \"77\u00d749=3773\"),\n    @(\"46\u00d711=506\", \"74\u00d765=4810\"),\n    @(\"36\u00d751=1836\", \"64\u00d757=3648\"),\n    @(\"37\u00d766=2442\", \"77\u00d776=5852\"),\n    @(\"94\u00d729=2726\", \"67\u00d761=4087\"),\n    @(\"67\u00d715=1005\", \"17\u00d719=323\"),\n    @(\"89\u00d762=5518\", \"92\u00d720=1840\"),\n    @(\"63\u00d754=3402\", \"59\u00d752=3068\"),\n    @(\"59\u00d720=1180\", \"25\u00d713=325\"),\n    @(\"87\u00d740=3480\", \"22\u00d724=528\"),\n    @(\"32\u00d787=2784\", \"92\u00d725=2300\"),\n    @(\"60\u00d734=2040\", \"54\u00d772=3888\"),\n    @(\"69\u00d794=6486\", \"16\u00d780=1280\"),\n    @(\"78\u00d719=1482\", \"20\u00d782=1640\"),\n    @(\"84\u00d783=6972\", \"65\u00d734=2210\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
